# Atualização de bases das ligas, do dia: 29-03-2024 às 13:24
# Bulgaria First League: swap the "Slavia Sofia" / "Botev Plovdiv" labels back
# to their correct fixtures, fix the swapped match-id pair (row 8/9), refresh
# the odds for the most recent fixtures (rows 213-218) and drop the two
# now-redundant trailing placeholder rows (219-220).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows that no longer exist after the refresh
# (this also shrinks the sheet's used range from AC220 down to AC218).
$ws.Range("A219:AC220").EntireRow.Delete()

$ws.Range("G2").Value2 = "Slavia Sofia"
$ws.Range("G3").Value2 = "Botev Plovdiv"
$ws.Range("AA8").Value2 = 1.025
$ws.Range("AC8").Value2 = 0.425
$ws.Range("B8").Value2 = 6627736
$ws.Range("F8").Value2 = "Botev Plovdiv"
$ws.Range("G8").Value2 = "Arda Kardzhali"
$ws.Range("H8").Value2 = 0
$ws.Range("I8").Value2 = 3
$ws.Range("J8").Value2 = "A"
$ws.Range("K8").Value2 = 5.25
$ws.Range("L8").Value2 = 3.6
$ws.Range("M8").Value2 = 1.571
$ws.Range("N8").Value2 = 26
$ws.Range("O8").Value2 = 11
$ws.Range("P8").Value2 = 1.083
$ws.Range("Q8").Value2 = 2.5
$ws.Range("R8").Value2 = 1.825
$ws.Range("S8").Value2 = 2.025
$ws.Range("T8").Value2 = 3.25
$ws.Range("U8").Value2 = 2
$ws.Range("V8").Value2 = 1.85
$ws.Range("W8").Value2 = -1
$ws.Range("Y8").Value2 = 0.08299999999999996
$ws.Range("Z8").Value2 = -1
$ws.Range("AA9").Value2 = -1
$ws.Range("AC9").Value2 = 0.4875
$ws.Range("B9").Value2 = 6627737
$ws.Range("F9").Value2 = "Slavia Sofia"
$ws.Range("G9").Value2 = "Lokomotiv 1929 Sofia"
$ws.Range("H9").Value2 = 2
$ws.Range("I9").Value2 = 0
$ws.Range("J9").Value2 = "H"
$ws.Range("K9").Value2 = 1.5
$ws.Range("L9").Value2 = 3.75
$ws.Range("M9").Value2 = 6.5
$ws.Range("N9").Value2 = 1.444
$ws.Range("O9").Value2 = 4.333
$ws.Range("P9").Value2 = 8
$ws.Range("Q9").Value2 = -1.25
$ws.Range("R9").Value2 = 2
$ws.Range("S9").Value2 = 1.85
$ws.Range("T9").Value2 = 2.25
$ws.Range("U9").Value2 = 1.875
$ws.Range("V9").Value2 = 1.975
$ws.Range("W9").Value2 = 0.444
$ws.Range("Y9").Value2 = -1
$ws.Range("Z9").Value2 = 1
$ws.Range("G18").Value2 = "Slavia Sofia"
$ws.Range("G21").Value2 = "Botev Plovdiv"
$ws.Range("F24").Value2 = "Slavia Sofia"
$ws.Range("G28").Value2 = "Botev Plovdiv"
$ws.Range("G32").Value2 = "Slavia Sofia"
$ws.Range("F36").Value2 = "Slavia Sofia"
$ws.Range("G39").Value2 = "Botev Plovdiv"
$ws.Range("G45").Value2 = "Slavia Sofia"
$ws.Range("F46").Value2 = "Botev Plovdiv"
$ws.Range("G53").Value2 = "Botev Plovdiv"
$ws.Range("F55").Value2 = "Slavia Sofia"
$ws.Range("G62").Value2 = "Slavia Sofia"
$ws.Range("F64").Value2 = "Botev Plovdiv"
$ws.Range("F67").Value2 = "Slavia Sofia"
$ws.Range("G67").Value2 = "Botev Plovdiv"
$ws.Range("F75").Value2 = "Botev Plovdiv"
$ws.Range("G76").Value2 = "Slavia Sofia"
$ws.Range("F83").Value2 = "Slavia Sofia"
$ws.Range("G85").Value2 = "Botev Plovdiv"
$ws.Range("G90").Value2 = "Slavia Sofia"
$ws.Range("F92").Value2 = "Botev Plovdiv"
$ws.Range("F98").Value2 = "Slavia Sofia"
$ws.Range("G103").Value2 = "Botev Plovdiv"
$ws.Range("F108").Value2 = "Botev Plovdiv"
$ws.Range("G111").Value2 = "Slavia Sofia"
$ws.Range("G115").Value2 = "Slavia Sofia"
$ws.Range("G118").Value2 = "Botev Plovdiv"
$ws.Range("F122").Value2 = "Slavia Sofia"
$ws.Range("F125").Value2 = "Botev Plovdiv"
$ws.Range("F131").Value2 = "Botev Plovdiv"
$ws.Range("F133").Value2 = "Slavia Sofia"
$ws.Range("F136").Value2 = "Botev Plovdiv"
$ws.Range("F140").Value2 = "Botev Plovdiv"
$ws.Range("G146").Value2 = "Slavia Sofia"
$ws.Range("F150").Value2 = "Slavia Sofia"
$ws.Range("G154").Value2 = "Botev Plovdiv"
$ws.Range("G156").Value2 = "Slavia Sofia"
$ws.Range("F157").Value2 = "Botev Plovdiv"
$ws.Range("G165").Value2 = "Botev Plovdiv"
$ws.Range("F168").Value2 = "Slavia Sofia"
$ws.Range("G174").Value2 = "Slavia Sofia"
$ws.Range("F175").Value2 = "Botev Plovdiv"
$ws.Range("F182").Value2 = "Slavia Sofia"
$ws.Range("G187").Value2 = "Botev Plovdiv"
$ws.Range("F192").Value2 = "Botev Plovdiv"
$ws.Range("G192").Value2 = "Slavia Sofia"
$ws.Range("G198").Value2 = "Botev Plovdiv"
$ws.Range("F203").Value2 = "Slavia Sofia"
$ws.Range("F207").Value2 = "Botev Plovdiv"
$ws.Range("G212").Value2 = "Slavia Sofia"
$ws.Range("B213").Value2 = 7973584
$ws.Range("E213").Value2 = 45380.58333333334
$ws.Range("F213").Value2 = "Levski Sofia"
$ws.Range("G213").Value2 = "Botev Plovdiv"
$ws.Range("K213").Value2 = 1.444
$ws.Range("L213").Value2 = 4.2
$ws.Range("M213").Value2 = 7.5
$ws.Range("N213").Value2 = 1.533
$ws.Range("O213").Value2 = 3.8
$ws.Range("P213").Value2 = 7
$ws.Range("Q213").Value2 = -1
$ws.Range("U213").Value2 = 1.9
$ws.Range("V213").Value2 = 1.95
$ws.Range("B214").Value2 = 7973585
$ws.Range("E214").Value2 = 45381.375
$ws.Range("F214").Value2 = "Botev Vratsa"
$ws.Range("G214").Value2 = "FC Hebar Pazardzhik"
$ws.Range("K214").Value2 = 2.4
$ws.Range("L214").Value2 = 3.2
$ws.Range("M214").Value2 = 3
$ws.Range("N214").Value2 = 2.375
$ws.Range("O214").Value2 = 3.25
$ws.Range("P214").Value2 = 3.1
$ws.Range("Q214").Value2 = -0.25
$ws.Range("R214").Value2 = 2.05
$ws.Range("S214").Value2 = 1.8
$ws.Range("U214").Value2 = 1.8
$ws.Range("V214").Value2 = 2.05
$ws.Range("B215").Value2 = 7973586
$ws.Range("E215").Value2 = 45381.47916666666
$ws.Range("F215").Value2 = "Lokomotiv Plovdiv"
$ws.Range("G215").Value2 = "Ludogorets Razgrad"
$ws.Range("K215").Value2 = 5.25
$ws.Range("L215").Value2 = 3.75
$ws.Range("M215").Value2 = 1.65
$ws.Range("N215").Value2 = 7.5
$ws.Range("O215").Value2 = 4.2
$ws.Range("P215").Value2 = 1.45
$ws.Range("Q215").Value2 = 1.25
$ws.Range("R215").Value2 = 1.8
$ws.Range("S215").Value2 = 2.05
$ws.Range("T215").Value2 = 2.75
$ws.Range("U215").Value2 = 2.025
$ws.Range("V215").Value2 = 1.825
$ws.Range("B216").Value2 = 6978434
$ws.Range("E216").Value2 = 45381.58333333334
$ws.Range("F216").Value2 = "Lokomotiv 1929 Sofia"
$ws.Range("G216").Value2 = "CSKA Sofia"
$ws.Range("K216").Value2 = 7.5
$ws.Range("L216").Value2 = 5
$ws.Range("M216").Value2 = 1.363
$ws.Range("N216").Value2 = 19
$ws.Range("O216").Value2 = 5.5
$ws.Range("P216").Value2 = 1.222
$ws.Range("Q216").Value2 = 1.75
$ws.Range("R216").Value2 = 1.85
$ws.Range("S216").Value2 = 2
$ws.Range("T216").Value2 = 2.5
$ws.Range("U216").Value2 = 1.975
$ws.Range("V216").Value2 = 1.875
$ws.Range("B217").Value2 = 7973587
$ws.Range("E217").Value2 = 45382.4375
$ws.Range("F217").Value2 = "Slavia Sofia"
$ws.Range("G217").Value2 = "Cherno More Varna"
$ws.Range("K217").Value2 = 2.5
$ws.Range("L217").Value2 = 3.1
$ws.Range("M217").Value2 = 2.9
$ws.Range("N217").Value2 = 2.2
$ws.Range("O217").Value2 = 3.2
$ws.Range("P217").Value2 = 3.5
$ws.Range("Q217").Value2 = -0.25
$ws.Range("R217").Value2 = 1.975
$ws.Range("S217").Value2 = 1.875
$ws.Range("T217").Value2 = 2.25
$ws.Range("U217").Value2 = 1.925
$ws.Range("V217").Value2 = 1.925
$ws.Range("B218").Value2 = 7973588
$ws.Range("E218").Value2 = 45382.54166666666
$ws.Range("F218").Value2 = "CSKA 1948 Sofia"
$ws.Range("G218").Value2 = "Beroe"
$ws.Range("K218").Value2 = 1.615
$ws.Range("L218").Value2 = 3.6
$ws.Range("M218").Value2 = 6
$ws.Range("N218").Value2 = 1.6
$ws.Range("O218").Value2 = 3.6
$ws.Range("P218").Value2 = 6
$ws.Range("Q218").Value2 = -1
$ws.Range("R218").Value2 = 2.05
$ws.Range("S218").Value2 = 1.8
$ws.Range("T218").Value2 = 2.25
$ws.Range("U218").Value2 = 1.8
$ws.Range("V218").Value2 = 2.05
